# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (rows 16-36, columns C:G) is
# re-sorted by "Periodo Mora" (column E): the new period 1712 record for
# MARIA EUGENIA DURANGO PINTO (previously the very last row, with a
# reduced "Valor Mora" of 14754) now leads the table, followed by every
# worker's 1805 record, and finally every worker's 1806 record.
# The underlying (worker, period, valor mora, salario basico) facts are
# unchanged - only their row order is rewritten here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    @{ Row = 16; C = "1050969804";  D = "MARIA EUGENIA DURANGO PINTO";  E = "1712"; F = 14754; G = 781242  },
    @{ Row = 17; C = "45486132";    D = "MARIA EUGENIA HERAZO CONEO";   E = "1805"; F = 23437; G = 781242  },
    @{ Row = 18; C = "73162110";    D = "JUAN ANTONIO RAMOS MERCADO";  E = "1805"; F = 23437; G = 781242  },
    @{ Row = 19; C = "1051818784";  D = "JUAN DARIO LOMBANA HERRERA";   E = "1805"; F = 23437; G = 781242  },
    @{ Row = 20; C = "45507230";    D = "MABEL KARINA ANGULO NOVOA";    E = "1805"; F = 23437; G = 781242  },
    @{ Row = 21; C = "45760350";    D = "JOANA ELAINE TAWIL DOMINGUEZ"; E = "1805"; F = 45000; G = 1500000 },
    @{ Row = 22; C = "9239499";     D = "OSMIN ENRIQUE ORTEGA ARROYO";  E = "1805"; F = 23437; G = 781242  },
    @{ Row = 23; C = "73129670";    D = "FELIX VASQUEZ ACEVEDO";        E = "1805"; F = 23437; G = 781242  },
    @{ Row = 24; C = "1143325736";  D = "DEIBER EDUARDO DIAGO BUELVAS"; E = "1805"; F = 23437; G = 781242  },
    @{ Row = 25; C = "15249140";    D = "WALBERTO ENRIQUE PAJARO DIAZ"; E = "1805"; F = 23437; G = 781242  },
    @{ Row = 26; C = "1050969804";  D = "MARIA EUGENIA DURANGO PINTO";  E = "1805"; F = 23437; G = 781242  },
    @{ Row = 27; C = "45486132";    D = "MARIA EUGENIA HERAZO CONEO";   E = "1806"; F = 23437; G = 781242  },
    @{ Row = 28; C = "73162110";    D = "JUAN ANTONIO RAMOS MERCADO";  E = "1806"; F = 23437; G = 781242  },
    @{ Row = 29; C = "1051818784";  D = "JUAN DARIO LOMBANA HERRERA";   E = "1806"; F = 23437; G = 781242  },
    @{ Row = 30; C = "45507230";    D = "MABEL KARINA ANGULO NOVOA";    E = "1806"; F = 23437; G = 781242  },
    @{ Row = 31; C = "45760350";    D = "JOANA ELAINE TAWIL DOMINGUEZ"; E = "1806"; F = 45000; G = 1500000 },
    @{ Row = 32; C = "9239499";     D = "OSMIN ENRIQUE ORTEGA ARROYO";  E = "1806"; F = 23437; G = 781242  },
    @{ Row = 33; C = "73129670";    D = "FELIX VASQUEZ ACEVEDO";        E = "1806"; F = 23437; G = 781242  },
    @{ Row = 34; C = "1143325736";  D = "DEIBER EDUARDO DIAGO BUELVAS"; E = "1806"; F = 23437; G = 781242  },
    @{ Row = 35; C = "15249140";    D = "WALBERTO ENRIQUE PAJARO DIAZ"; E = "1806"; F = 23437; G = 781242  },
    @{ Row = 36; C = "1050969804";  D = "MARIA EUGENIA DURANGO PINTO";  E = "1806"; F = 23437; G = 781242  }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
}

"Estado de cuenta rows refreshed"
